# Update NATMI LR-pair sheet (Efnb3-Ephb4) with refreshed TPM-based results.
# Old data had 9 sending/target cluster combinations (ECs/FAPs/MuSCs x3),
# spread across rows 2-10. The refreshed run only has 6 rows of data
# (rows 2-7); the former "MuSCs" sending-cluster block (old rows 8-10) is
# gone and the remaining rows carry updated metric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old rows 8-10 (sheet shrinks from A1:T10 to A1:T7).
$ws.Range("A8:A10").EntireRow.Delete()

# Row 2: FAPs -> Efnb3 -> Ephb4 -> ECs
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Efnb3"
$ws.Cells.Item(2, 3).Value = "Ephb4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.1498043333333333
$ws.Cells.Item(2, 8).Value = 0.449413
$ws.Cells.Item(2, 9).Value = 0.08722868471333377
$ws.Cells.Item(2, 10).Value = 0.08722868471333377
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 19.11595033333333
$ws.Cells.Item(2, 14).Value = 57.347851
$ws.Cells.Item(2, 15).Value = 0.6851940154453416
$ws.Cells.Item(2, 16).Value = 0.6851940154453418
$ws.Cells.Item(2, 17).Value = 2.863652195718112
$ws.Cells.Item(2, 18).Value = 25.772869761463
$ws.Cells.Item(2, 19).Value = 0.05976857274074485
$ws.Cells.Item(2, 20).Value = 0.05976857274074487

# Row 3: FAPs -> Efnb3 -> Ephb4 -> FAPs
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Efnb3"
$ws.Cells.Item(3, 3).Value = "Ephb4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.1498043333333333
$ws.Cells.Item(3, 8).Value = 0.449413
$ws.Cells.Item(3, 9).Value = 0.08722868471333377
$ws.Cells.Item(3, 10).Value = 0.08722868471333377
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 4.865208333333334
$ws.Cells.Item(3, 14).Value = 14.595625
$ws.Cells.Item(3, 15).Value = 0.1743890089566637
$ws.Cells.Item(3, 16).Value = 0.1743890089566637
$ws.Cells.Item(3, 17).Value = 0.7288292909027779
$ws.Cells.Item(3, 18).Value = 6.559463618125001
$ws.Cells.Item(3, 19).Value = 0.01521172387975156
$ws.Cells.Item(3, 20).Value = 0.01521172387975156

# Row 4: FAPs -> Efnb3 -> Ephb4 -> MuSCs
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Efnb3"
$ws.Cells.Item(4, 3).Value = "Ephb4"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.1498043333333333
$ws.Cells.Item(4, 8).Value = 0.449413
$ws.Cells.Item(4, 9).Value = 0.08722868471333377
$ws.Cells.Item(4, 10).Value = 0.08722868471333377
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 3.917436333333333
$ws.Cells.Item(4, 14).Value = 11.752309
$ws.Cells.Item(4, 15).Value = 0.1404169755979945
$ws.Cells.Item(4, 16).Value = 0.1404169755979946
$ws.Cells.Item(4, 17).Value = 0.5868489382907778
$ws.Cells.Item(4, 18).Value = 5.281640444617
$ws.Cells.Item(4, 19).Value = 0.01224838809283735
$ws.Cells.Item(4, 20).Value = 0.01224838809283735

# Row 5: MuSCs -> Efnb3 -> Ephb4 -> ECs
$ws.Cells.Item(5, 1).Value = "MuSCs"
$ws.Cells.Item(5, 2).Value = "Efnb3"
$ws.Cells.Item(5, 3).Value = "Ephb4"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.567570333333333
$ws.Cells.Item(5, 8).Value = 4.702711
$ws.Cells.Item(5, 9).Value = 0.9127713152866662
$ws.Cells.Item(5, 10).Value = 0.9127713152866662
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 19.11595033333333
$ws.Cells.Item(5, 14).Value = 57.347851
$ws.Cells.Item(5, 15).Value = 0.6851940154453416
$ws.Cells.Item(5, 16).Value = 0.6851940154453418
$ws.Cells.Item(5, 17).Value = 29.96559663600678
$ws.Cells.Item(5, 18).Value = 269.690369724061
$ws.Cells.Item(5, 19).Value = 0.6254254427045968
$ws.Cells.Item(5, 20).Value = 0.6254254427045969

# Row 6: MuSCs -> Efnb3 -> Ephb4 -> FAPs
$ws.Cells.Item(6, 1).Value = "MuSCs"
$ws.Cells.Item(6, 2).Value = "Efnb3"
$ws.Cells.Item(6, 3).Value = "Ephb4"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.567570333333333
$ws.Cells.Item(6, 8).Value = 4.702711
$ws.Cells.Item(6, 9).Value = 0.9127713152866662
$ws.Cells.Item(6, 10).Value = 0.9127713152866662
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 4.865208333333334
$ws.Cells.Item(6, 14).Value = 14.595625
$ws.Cells.Item(6, 15).Value = 0.1743890089566637
$ws.Cells.Item(6, 16).Value = 0.1743890089566637
$ws.Cells.Item(6, 17).Value = 7.626556248819446
$ws.Cells.Item(6, 18).Value = 68.63900623937501
$ws.Cells.Item(6, 19).Value = 0.1591772850769121
$ws.Cells.Item(6, 20).Value = 0.1591772850769122

# Row 7: MuSCs -> Efnb3 -> Ephb4 -> MuSCs
$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Efnb3"
$ws.Cells.Item(7, 3).Value = "Ephb4"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.567570333333333
$ws.Cells.Item(7, 8).Value = 4.702711
$ws.Cells.Item(7, 9).Value = 0.9127713152866662
$ws.Cells.Item(7, 10).Value = 0.9127713152866662
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.917436333333333
$ws.Cells.Item(7, 14).Value = 11.752309
$ws.Cells.Item(7, 15).Value = 0.1404169755979945
$ws.Cells.Item(7, 16).Value = 0.1404169755979946
$ws.Cells.Item(7, 17).Value = 6.140856978855444
$ws.Cells.Item(7, 18).Value = 55.267712809699
$ws.Cells.Item(7, 19).Value = 0.1281685875051572
$ws.Cells.Item(7, 20).Value = 0.1281685875051572
